# Merge the "xlong" / "ylat" columns into a single "geopoint" column
# (Table Schema's geopoint type), dropping the now redundant ylat column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Rename the header of column D from "xlong" to "geopoint"
$ws.Range("D1").Value = "geopoint"

# Combine "xlong" (col D) and "ylat" (col E) into "xlong, ylat" in col D,
# for every data row, before column E disappears.
for ($r = 2; $r -le $lastRow; $r++) {
    $xlong = $ws.Cells.Item($r, 4).Value()
    $ylat = $ws.Cells.Item($r, 5).Value()
    if (($xlong -ne $null) -and ($ylat -ne $null)) {
        $ws.Cells.Item($r, 4).Value = "$xlong, $ylat"
    }
}

# Remove the now-redundant ylat column; adresse/emplacements/no_appel/info
# shift one column to the left (F->E, G->F, H->G, I->H).
$ws.Range("E:E").Delete() | Out-Null

# The phone number (now in column G, row 3) should be stored as a real number.
$ws.Cells.Item(3, 7).Value = 476544254

# Widen the merged geopoint column so the combined values fit.
$ws.Range("D1").ColumnWidth = 19.35

# Reset the active selection to A1.
$ws.Range("A1").Select() | Out-Null
